$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.877
$ws.Range("C8").Value = -11.953
$ws.Range("C10").Value = -13.08
$ws.Range("C12").Value = -11.466
$ws.Range("D12").Value = -7.775
$ws.Range("D15").Value = -8.011999999999999
$ws.Range("D17").Value = -8.228999999999999
$ws.Range("C18").Value = -12.31
$ws.Range("D26").Value = -7.255000000000001
$ws.Range("D27").Value = -7.572
$ws.Range("D28").Value = -7.609
$ws.Range("C37").Value = -13.179
$ws.Range("D37").Value = -7.369999999999999
$ws.Range("D47").Value = -7.452
$ws.Range("C55").Value = -13.633
$ws.Range("D65").Value = -7.741
$ws.Range("C68").Value = -10.809
$ws.Range("D73").Value = -7.826000000000001
$ws.Range("C77").Value = -13.169
$ws.Range("C78").Value = -13.004
$ws.Range("C81").Value = -13.054
$ws.Range("C82").Value = -11.773
$ws.Range("D84").Value = -7.975
$ws.Range("D85").Value = -8.632
$ws.Range("D93").Value = -7.131
$ws.Range("D95").Value = -7.506
$ws.Range("D98").Value = -7.231
$ws.Range("D99").Value = -8.193000000000001
$ws.Range("D101").Value = -7.825
